# Updates the cryptocurrency price/volume data in the active worksheet
# to reflect the latest scrape (GitHub Actions scheduled update).
# Numeric-looking price strings are written via a temporary Text
# number format so Excel keeps them as literal strings (e.g. "61.00")
# instead of auto-converting them to numbers; the style is then reset
# back to Normal so no extra cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.369.72"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.606.76"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "1.839.31"
$ws.Range("D13").Value = "1.636.97"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.505"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "26.374.88"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.52%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.17%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "1.445.72"
$ws.Range("E33").Value = "  +8.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0165"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.827"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.938"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("D44").Value = "1.746.85"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("E49").Value = "  -2.55%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0952"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.17%  "